$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing employee's hourly_rate (C2: 20 -> 80)
$ws.Range("C2").Value = 80

# Add new employee rows
# Row 3: Tee
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Tee"
$ws.Range("C3").Value = 5
$ws.Range("D3").Value = 7
$ws.Range("E3").Value = 9

# Row 4: CSL
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "CSL"
$ws.Range("C4").Value = 9
$ws.Range("D4").Value = 8
$ws.Range("E4").Value = 10
